$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/22/2024  Through  7/28/2024"

# --- Row 15 ---
$ws.Range("A14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("N15").Value = -70.588235294117

# --- Row 16 ---
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -44.444444444444
$ws.Range("I16").Value = 72
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = 16.129032258064
$ws.Range("L16").Value = 10.76923076923
$ws.Range("M16").Value = -46.268656716417
$ws.Range("N16").Value = -81.679389312977

# --- Row 17 ---
$ws.Range("A14").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = "0"
$ws.Range("E17").Value = "***.*"
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 70
$ws.Range("I17").Value = 94
$ws.Range("K17").Value = -3.092783505154
$ws.Range("L17").Value = -2.083333333333
$ws.Range("M17").Value = 36.231884057971
$ws.Range("N17").Value = -47.777777777777

# --- Row 18 ---
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 38
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = -28.301886792452
$ws.Range("L18").Value = -38.709677419354
$ws.Range("M18").Value = -73.793103448275
$ws.Range("N18").Value = -94.516594516594

# --- Row 19 ---
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = -55.555555555555
$ws.Range("I19").Value = 323
$ws.Range("J19").Value = 370
$ws.Range("K19").Value = -12.702702702702
$ws.Range("L19").Value = 13.333333333333
$ws.Range("M19").Value = 22.348484848484
$ws.Range("N19").Value = -15.44502617801

# --- Row 20 ---
$ws.Range("C20").Value = 8
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 18
$ws.Range("H20").Value = -10
$ws.Range("I20").Value = 95
$ws.Range("J20").Value = 73
$ws.Range("K20").Value = 30.136986301369
$ws.Range("L20").Value = 26.666666666666
$ws.Range("M20").Value = -2.061855670103
$ws.Range("N20").Value = -94.099378881987

# --- Row 21 ---
$ws.Range("C21").Value = 25
$ws.Range("E21").Value = -13.793103448275
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 116
$ws.Range("H21").Value = -31.896551724137
$ws.Range("I21").Value = 627
$ws.Range("J21").Value = 664
$ws.Range("K21").Value = -5.572289156626
$ws.Range("L21").Value = 6.271186440677
$ws.Range("M21").Value = -12.67409470752
$ws.Range("N21").Value = -80.919050517346

# --- Row 23 ---
$ws.Range("A14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("D23").Value = "0"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").Value = "0"
$ws.Range("H23").Value = -100

# --- Row 24 ---
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -8.333333333333
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = -7.54716981132
$ws.Range("I24").Value = 914
$ws.Range("J24").Value = 699
$ws.Range("K24").Value = 30.758226037196
$ws.Range("L24").Value = 44.849445324881
$ws.Range("M24").Value = 74.42748091603

# --- Row 25 ---
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 83
$ws.Range("G25").Value = 87
$ws.Range("H25").Value = -4.597701149425
$ws.Range("I25").Value = 754
$ws.Range("J25").Value = 489
$ws.Range("K25").Value = 54.192229038854
$ws.Range("L25").Value = 90.40404040404

# --- Row 26 ---
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 185
$ws.Range("J26").Value = 152
$ws.Range("K26").Value = 21.710526315789
$ws.Range("L26").Value = 43.410852713178
$ws.Range("M26").Value = -7.960199004975

# --- Row 27 ---
$ws.Range("J14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = -46.666666666666
$ws.Range("L27").Value = -33.333333333333

# --- Row 28 ---
$ws.Range("A14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 250
$ws.Range("I28").Value = 21
$ws.Range("K28").Value = 16.666666666666
$ws.Range("L28").Value = 5

$excel.Application.CutCopyMode = $false
